$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script re-ran and found nothing new for 2025-05-15, so the "NA" marker
# that used to sit in C111 moves down to the newly appended row.
# A bare quote-prefix forces an empty *text* cell (plain $null/"" would clear
# the cell entirely instead of leaving an empty inline string), then the
# style is reset to Normal so the quote-prefix doesn't leave a stray
# "stored as text" number format on the cell.
$ws.Range("C111").Value = "'"
$ws.Range("C111").Style = "Normal"

# New row 112: continuation of the "Rien ne nous concerne aujourd'hui !" series
# for 2025-05-16, ending in today's "NA" marker in column C.
$ws.Range("A112").Value = "'2025-05-16"
$ws.Range("A112").Style = "Normal"
$ws.Range("B112").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C112").Value = "NA"
$ws.Range("D112").Value = 1
